$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and B (student DNI / enrollment code) become Text-formatted
# for rows 2-25 (everything below the header row).
$ws.Range("A2:B25").NumberFormat = "@"

# New column L ("Ingresos") gets a 0 entered for every student row,
# centered both horizontally and vertically - matching the style already
# used for the grade columns (D:J) but with the default General number
# format instead of the "00" one.
$first = $ws.Range("L2")
$first.Value = 0
$first.HorizontalAlignment = -4108   # xlCenter
$first.VerticalAlignment = -4108     # xlCenter

$first.Copy()
$rest = $ws.Range("L3:L25")
$rest.PasteSpecial(-4122)            # xlPasteFormats
$rest.Value = 0

# Leave the new column selected, like in the saved workbook.
$ws.Range("L2:L25").Select() | Out-Null
